# filter function, layout tweaks, script fixes
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # cikkek
$ws2 = $wb.Worksheets.Item(2)   # vélemények

# --- Data edits on "cikkek" sheet ---
# E7 picks up the tag set that used to live on E5 ("web, branding, logo, frontend, case study")
$ws1.Range("E7").Value = $ws1.Range("E5").Value2
# E5 gets a new, more specific tag set (new shared string)
$ws1.Range("E5").Value = "web, branding, logo, uiux, case study"

# New "Kész" status cell at C11, copying the formatting used by the other
# "Kész" status cells (e.g. C2) so it gets the same fill style.
$ws1.Range("C2").Copy()
$ws1.Range("C11").PasteSpecial(-4122)  # xlPasteFormats
$ws1.Range("C11").Value = $ws1.Range("C2").Value2

# --- Layout / selection tweaks ---
# Move the active tab/selection from "vélemények" to "cikkek", and update
# the lingering selection left behind on "vélemények".
[void]$ws2.Range("B11").Select()
[void]$ws1.Activate()
[void]$ws1.Range("F4").Select()
